$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string; these must stay
# stored as TEXT (matching the source inlineStr cells), so temporarily force
# the Text number format before assigning, then clear the format again so the
# cell style reverts to the sheet default (General / no explicit style index).
$textCells = @('D5', 'D6', 'D7', 'D9', 'D10', 'D11', 'D12', 'D15', 'D19', 'D21', 'D23', 'D24', 'D27', 'D30', 'D31', 'D32', 'D34', 'D38', 'D40', 'D41', 'D42', 'D44', 'D45', 'D47', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '43.972.14'
$ws.Range("E2").Value = '  +3.31%  '
$ws.Range("D3").Value = '2.343.06'
$ws.Range("E3").Value = '  +2.66%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '312.66'
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("D6").Value = '109.03'
$ws.Range("D7").Value = '0.633'
$ws.Range("E7").Value = '  +1.13%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '0.620'
$ws.Range("E9").Value = '  +2.94%  '
$ws.Range("D10").Value = '41.29'
$ws.Range("E10").Value = '  +4.45%  '
$ws.Range("D11").Value = '0.0920'
$ws.Range("E11").Value = '  +2.02%  '
$ws.Range("D12").Value = '8.61'
$ws.Range("E12").Value = '  +2.66%  '
$ws.Range("E13").Value = '  +1.69%  '
$ws.Range("E14").Value = '  -1.00%  '
$ws.Range("D15").Value = '15.52'
$ws.Range("E15").Value = '  +2.20%  '
$ws.Range("D16").Value = '2.698.17'
$ws.Range("E16").Value = '  +2.61%  '
$ws.Range("D17").Value = '2.337.97'
$ws.Range("E17").Value = '  +1.27%  '
$ws.Range("D18").Value = '43.878.50'
$ws.Range("E18").Value = '  +2.38%  '
$ws.Range("D19").Value = '7.59'
$ws.Range("E19").Value = '  +2.93%  '
$ws.Range("E20").Value = '  +1.87%  '
$ws.Range("D21").Value = '13.00'
$ws.Range("E21").Value = '  -3.64%  '
$ws.Range("E22").Value = '  +0.65%  '
$ws.Range("D23").Value = '3.49'
$ws.Range("E23").Value = '  -0.99%  '
$ws.Range("D24").Value = '269.18'
$ws.Range("E24").Value = '  +1.97%  '
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("D27").Value = '7.58'
$ws.Range("E27").Value = '  +7.24%  '
$ws.Range("E28").Value = '  +3.32%  '
$ws.Range("E29").Value = '  -2.16%  '
$ws.Range("D30").Value = '38.83'
$ws.Range("E30").Value = '  +4.98%  '
$ws.Range("D31").Value = '22.72'
$ws.Range("E31").Value = '  +1.30%  '
$ws.Range("D32").Value = '168.31'
$ws.Range("E32").Value = '  +1.14%  '
$ws.Range("E33").Value = '  +1.94%  '
$ws.Range("D34").Value = '2.78'
$ws.Range("E34").Value = '  +7.59%  '
$ws.Range("E35").Value = '  +1.39%  '
$ws.Range("E36").Value = '  +5.28%  '
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("D38").Value = '0.0366'
$ws.Range("E38").Value = '  +4.61%  '
$ws.Range("E39").Value = '  +8.78%  '
$ws.Range("D40").Value = '3.82'
$ws.Range("E40").Value = '  +0.80%  '
$ws.Range("D41").Value = '1.71'
$ws.Range("E41").Value = '  +8.82%  '
$ws.Range("D42").Value = '104.76'
$ws.Range("E42").Value = '  +10.43%  '
$ws.Range("E43").Value = '  +3.15%  '
$ws.Range("D44").Value = '71.77'
$ws.Range("E44").Value = '  +2.14%  '
$ws.Range("D45").Value = '13.28'
$ws.Range("E45").Value = '  +9.23%  '
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").Value = '114.03'
$ws.Range("E47").Value = '  +0.39%  '
$ws.Range("D48").Value = '1.665.01'
$ws.Range("E48").Value = '  -3.65%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = '9.00'
$ws.Range("E49").Value = '  +3.06%  '
$ws.Range("B50").Value = 'ordi'
$ws.Range("C50").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D50").Value = '76.51'
$ws.Range("E50").Value = '  -3.20%  '
$ws.Range("D51").Value = '0.215'
$ws.Range("E51").Value = '  +13.62%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
